# Generate Report for Handback
# Updates the localization-status report after a handback: status moves
# from "Ready for handoff" to "Handed back: in sync with en-US", refreshed
# handback timestamps are recorded, and the (now resolved) handback-version
# warning is cleared from the per-locale "Error Detail" column.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet ------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus   # zh-cn status
$overview.Range("F2").Value = $newStatus   # de-de status

# ---- zh-cn sheet -----------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus                    # Status
$zhcn.Range("K2").Value = "2016-08-27 18:56:28"          # Latest Handback DateTime
$zhcn.Range("P2").Value = ""                             # Error Detail (resolved)

# ---- de-de sheet -----------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus                    # Status
$dede.Range("K2").Value = "2016-08-27 18:56:35"          # Latest Handback DateTime
$dede.Range("P2").Value = ""                             # Error Detail (resolved)

# ---- Column widths: reflow the Status / Error Detail columns -------------
$overview.Columns.Item(5).AutoFit()
$overview.Columns.Item(6).AutoFit()
$zhcn.Columns.Item(3).AutoFit()
$zhcn.Columns.Item(16).AutoFit()
$dede.Columns.Item(3).AutoFit()
$dede.Columns.Item(16).AutoFit()
